# Applies the per-row Price (D) / Volume(1h) (E) updates from the crypto-tracker
# refresh, plus the RenderToken/Dai row swap (rows 28-29), exactly as captured in
# the source diff. Values are written as literal text (matching the original
# inlineStr cells) - D-column entries that look numeric ("1.00", "10.25", ...) are
# written with a leading apostrophe so Excel stores them as text instead of
# coercing them to numbers (which would display "1" instead of "1.00"); the
# quote-prefix formatting Excel applies for that is then reset to Normal so the
# cell style stays untouched, matching the source file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '''67.564.42'
$ws.Range('E2').Value = '  -2.03%  '

# Row 3
$ws.Range('D3').Value = '''3.782.75'
$ws.Range('E3').Value = '  -0.05%  '

# Row 4
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  +0.01%  '

# Row 5
$ws.Range('D5').Value = '''593.66'
$ws.Range('E5').Value = '  -0.94%  '

# Row 6
$ws.Range('D6').Value = '''166.68'
$ws.Range('E6').Value = '  -1.90%  '

# Row 7
$ws.Range('D7').Value = '''3.782.94'
$ws.Range('E7').Value = '  +0.16%  '

# Row 8
$ws.Range('E8').Value = '  +0.02%  '

# Row 9
$ws.Range('D9').Value = '''0.518'
$ws.Range('E9').Value = '  -1.27%  '

# Row 10
$ws.Range('E10').Value = '  -1.95%  '

# Row 11
$ws.Range('E11').Value = '  -2.03%  '

# Row 12
$ws.Range('E12').Value = '  -1.06%  '

# Row 13
$ws.Range('E13').Value = '  -2.93%  '

# Row 14
$ws.Range('D14').Value = '''35.96'
$ws.Range('E14').Value = '  -2.14%  '

# Row 15
$ws.Range('D15').Value = '''4.418.77'
$ws.Range('E15').Value = '  -0.19%  '

# Row 16
$ws.Range('D16').Value = '''3.792.72'
$ws.Range('E16').Value = '  +0.22%  '

# Row 17
$ws.Range('D17').Value = '''67.491.58'
$ws.Range('E17').Value = '  -2.16%  '

# Row 18
$ws.Range('D18').Value = '''18.16'
$ws.Range('E18').Value = '  +0.11%  '

# Row 19
$ws.Range('E19').Value = '  +0.13%  '

# Row 20
$ws.Range('E20').Value = '  -1.21%  '

# Row 21
$ws.Range('D21').Value = '''10.25'
$ws.Range('E21').Value = '  -6.56%  '

# Row 22
$ws.Range('D22').Value = '''458.54'
$ws.Range('E22').Value = '  -2.65%  '

# Row 23
$ws.Range('D23').Value = '''0.697'
$ws.Range('E23').Value = '  -1.40%  '

# Row 24
$ws.Range('E24').Value = '  +2.85%  '

# Row 25
$ws.Range('D25').Value = '''83.53'
$ws.Range('E25').Value = '  -1.51%  '

# Row 26
$ws.Range('E26').Value = '  -4.11%  '

# Row 27
$ws.Range('E27').Value = '  -2.66%  '

# Row 28
$ws.Range('B28').Value = 'Dai'
$ws.Range('C28').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D28').Value = '''1.00'
$ws.Range('E28').Value = '  +0.08%  '

# Row 29
$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D29').Value = '''10.01'
$ws.Range('E29').Value = '  -2.29%  '

# Row 30
$ws.Range('D30').Value = '''2.78'
$ws.Range('E30').Value = '  -1.52%  '

# Row 31
$ws.Range('D31').Value = '''29.80'
$ws.Range('E31').Value = '  -1.50%  '

# Row 32
$ws.Range('E32').Value = '  -1.56%  '

# Row 33
$ws.Range('D33').Value = '''7.21'
$ws.Range('E33').Value = '  -3.18%  '

# Row 34
$ws.Range('D34').Value = '''9.16'
$ws.Range('E34').Value = '  -2.11%  '

# Row 35
$ws.Range('D35').Value = '''0.998'
$ws.Range('E35').Value = '  -0.05%  '

# Row 36
$ws.Range('D36').Value = '''3.736.32'
$ws.Range('E36').Value = '  -0.24%  '

# Row 37
$ws.Range('E37').Value = '  -2.21%  '

# Row 38
$ws.Range('E38').Value = '  -5.62%  '

# Row 39
$ws.Range('E39').Value = '  -1.09%  '

# Row 40
$ws.Range('D40').Value = '''0.994'
$ws.Range('E40').Value = '  -0.97%  '

# Row 41
$ws.Range('E41').Value = '  -2.22%  '

# Row 42
$ws.Range('D42').Value = '''1.00'
$ws.Range('E42').Value = '  +0.05%  '

# Row 43
$ws.Range('E43').Value = '  -0.01%  '

# Row 44
$ws.Range('D44').Value = '''44.04'
$ws.Range('E44').Value = '  +0.29%  '

# Row 45
$ws.Range('E45').Value = '  -3.89%  '

# Row 46
$ws.Range('D46').Value = '''47.01'
$ws.Range('E46').Value = '  +2.17%  '

# Row 47
$ws.Range('D47').Value = '''8.37'
$ws.Range('E47').Value = '  -3.19%  '

# Row 48
$ws.Range('D48').Value = '''147.48'
$ws.Range('E48').Value = '  +0.86%  '

# Row 49
$ws.Range('D49').Value = '''392.84'
$ws.Range('E49').Value = '  -2.04%  '

# Row 50
$ws.Range('E50').Value = '  -7.54%  '

# Row 51
$ws.Range('D51').Value = '''2.756.05'
$ws.Range('E51').Value = '  +2.21%  '

# Strip the auto-applied quote-prefix style from the numeric-looking text cells
# above so their formatting matches the source workbook (no explicit style).
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').Style = 'Normal'
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'

